$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Tutor "Environment:" line (paragraph ~17): add Visual Studio Code
#    "Git, Github and Heroku Deployment" -> "Visual Studio Code, Github and Heroku Deployment"
# ------------------------------------------------------------------
$p = $d.Paragraphs.Item(17)
$p.Range.Find.Execute("Git, Github and Heroku Deployment", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Visual Studio Code, Github and Heroku Deployment", 2) | Out-Null

# ------------------------------------------------------------------
# 2) Student "Environment:" line (paragraph ~22): same tool swap, and
#    "Photoshop" -> "Zoom"
#    "Git, Github and Heroku Deployment, Photoshop " -> "Visual Studio Code, Github and Heroku Deployment, Zoom"
# ------------------------------------------------------------------
$p = $d.Paragraphs.Item(22)
$p.Range.Find.Execute("Git, Github and Heroku Deployment, Photoshop ", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Visual Studio Code, Github and Heroku Deployment, Zoom", 2) | Out-Null

# ------------------------------------------------------------------
# 3) "Secured disability benefits and" -> "Secure disability benefits and"
# ------------------------------------------------------------------
$d.Content.Find.Execute("Secured disability benefits and", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Secure disability benefits and", 2) | Out-Null

# ------------------------------------------------------------------
# 4) Drop "Purple Finch Realty LLC" co-employer mention
#    "Sherrie McNulty LLC and Purple Finch Realty LLC in Merrimack, NH " -> "Sherrie McNulty LLC in Merrimack, NH "
# ------------------------------------------------------------------
$d.Content.Find.Execute("Sherrie McNulty LLC and Purple Finch Realty LLC in Merrimack, NH ", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Sherrie McNulty LLC in Merrimack, NH ", 2) | Out-Null

# ------------------------------------------------------------------
# 5) Simplify job title
#    "Computer Consultant/Owner Broker  " -> "Consultant "
# ------------------------------------------------------------------
$d.Content.Find.Execute("Computer Consultant/Owner Broker  ", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Consultant ", 2) | Out-Null

# ------------------------------------------------------------------
# 6) "Establish and maintain contacts" -> "Establish and maintain customer relationships"
# ------------------------------------------------------------------
$d.Content.Find.Execute("Establish and maintain contacts", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Establish and maintain customer relationships", 2) | Out-Null

# ------------------------------------------------------------------
# 7) "Determine customer needs" -> "Customize existing websites"
# ------------------------------------------------------------------
$d.Content.Find.Execute("Determine customer needs", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Customize existing websites", 2) | Out-Null

# ------------------------------------------------------------------
# 8) "Environment: HTML, CSS" (Sherrie McNulty LLC role) gains ", Windows"
#    Target only the bare "Environment: HTML, CSS" paragraph (not the
#    later "Environment: HTML, CSS, JavaScript, ..." one), so scope the
#    Find to that specific paragraph.
# ------------------------------------------------------------------
$p = $d.Paragraphs.Item(34)
$p.Range.Find.Execute("Environment: HTML, CSS", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Environment: HTML, CSS, Windows", 2) | Out-Null

# ------------------------------------------------------------------
# 9) Remove the two bullet paragraphs that are no longer relevant:
#    "Design, implement and install applications" and
#    "Systems administration"
# ------------------------------------------------------------------
$p = $d.Paragraphs.Item(32)
$p.Range.Delete() | Out-Null
$p = $d.Paragraphs.Item(32)
$p.Range.Delete() | Out-Null

# ------------------------------------------------------------------
# 10) Software Engineer bullets get reworked
#    "Design and develop web and internal applications" ->
#       "Implement web-based application to enable journalists to
#        publish and automatically move obsolete articles to a user
#        accessible archive without assistance from software developers"
# ------------------------------------------------------------------
$d.Content.Find.Execute("Design and develop web and internal applications", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Implement web-based application to enable journalists to publish and automatically move obsolete articles to a user accessible archive without assistance from software developers", 2) | Out-Null

# ------------------------------------------------------------------
# 11) "Customer support" -> "Generate technical documentation and user guide"
# ------------------------------------------------------------------
$d.Content.Find.Execute("Customer support", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Generate technical documentation and user guide", 2) | Out-Null

# ------------------------------------------------------------------
# 12) Former "Design and develop a web-based application..." bullet ->
#     "Customer training and support"
# ------------------------------------------------------------------
$d.Content.Find.Execute("Design and develop a web-based application that enabled journalists to publish and automatically archive articles without assistance from software developers", `
    $true, $false, $false, $false, $false, $true, 1, $false, "Customer training and support", 2) | Out-Null

# ------------------------------------------------------------------
# 13) Remove the two trailing bullet paragraphs:
#     "Document detailed design specifications and users guide" and
#     "Installation, training and customer support"
#     (originally paragraphs 41 & 42, shifted down by 2 because of the
#     earlier deletion in step 9)
# ------------------------------------------------------------------
$p = $d.Paragraphs.Item(39)
$p.Range.Delete() | Out-Null
$p = $d.Paragraphs.Item(39)
$p.Range.Delete() | Out-Null

# ------------------------------------------------------------------
# 14) "Shell Script, C" -> "Shell Script, UNIX, Mac, Windows"
# ------------------------------------------------------------------
$d.Content.Find.Execute("Shell Script, C", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Shell Script, UNIX, Mac, Windows", 2) | Out-Null

Write-Output "done"
